# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the crypto
# symbol list with the latest scraped values. Each value is written as
# literal text (leading apostrophe) so Excel doesn't coerce the price /
# percentage strings into numeric cells, then ClearFormats() drops the
# transient "stored as text" quote-prefix style so the cell format stays
# the same as before the write.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "318.83" },
    @{ Cell = "E2"; Value = "3.04%" },
    @{ Cell = "D3"; Value = "41.30" },
    @{ Cell = "E3"; Value = "0.87%" },
    @{ Cell = "D4"; Value = "5.237" },
    @{ Cell = "E4"; Value = "2.25%" },
    @{ Cell = "D5"; Value = "0.07703" },
    @{ Cell = "E5"; Value = "1.01%" },
    @{ Cell = "D6"; Value = "1.691" },
    @{ Cell = "E6"; Value = "5.14%" },
    @{ Cell = "D7"; Value = "0.9411" },
    @{ Cell = "E7"; Value = "3.42%" },
    @{ Cell = "D9"; Value = "0.1262" },
    @{ Cell = "E9"; Value = "-0.57%" },
    @{ Cell = "D10"; Value = "0.1847" },
    @{ Cell = "E10"; Value = "2.20%" },
    @{ Cell = "D11"; Value = "0.09225" },
    @{ Cell = "E11"; Value = "1.03%" },
    @{ Cell = "D12"; Value = "0.04127" },
    @{ Cell = "E12"; Value = "-4.77%" },
    @{ Cell = "E13"; Value = "0.30%" },
    @{ Cell = "E14"; Value = "3.29%" },
    @{ Cell = "D15"; Value = "0.006000" },
    @{ Cell = "E15"; Value = "4.83%" },
    @{ Cell = "D17"; Value = "3.350" },
    @{ Cell = "E17"; Value = "0.00%" },
    @{ Cell = "D18"; Value = "4.349" },
    @{ Cell = "E18"; Value = "1.49%" },
    @{ Cell = "D19"; Value = "0.3351" },
    @{ Cell = "E19"; Value = "1.01%" },
    @{ Cell = "D20"; Value = "8.403" },
    @{ Cell = "E20"; Value = "21.07%" },
    @{ Cell = "D21"; Value = "0.1357" },
    @{ Cell = "E21"; Value = "-2.57%" },
    @{ Cell = "D22"; Value = "0.2738" },
    @{ Cell = "E22"; Value = "-0.03%" },
    @{ Cell = "D23"; Value = "0.04029" },
    @{ Cell = "E23"; Value = "-0.40%" },
    @{ Cell = "D24"; Value = "0.001273" },
    @{ Cell = "E24"; Value = "0.21%" },
    @{ Cell = "D25"; Value = "0.004103" },
    @{ Cell = "E25"; Value = "0.13%" },
    @{ Cell = "D26"; Value = "0.0001276" },
    @{ Cell = "E26"; Value = "0.34%" },
    @{ Cell = "D38"; Value = "0.02524" },
    @{ Cell = "E38"; Value = "4.24%" },
    @{ Cell = "D39"; Value = "0.05328" },
    @{ Cell = "E39"; Value = "1.74%" },
    @{ Cell = "D40"; Value = "0.007757" },
    @{ Cell = "E40"; Value = "-1.27%" },
    @{ Cell = "D41"; Value = "0.1314" },
    @{ Cell = "E41"; Value = "0.94%" },
    @{ Cell = "D42"; Value = "0.007045" },
    @{ Cell = "E42"; Value = "3.47%" },
    @{ Cell = "D43"; Value = "0.002159" },
    @{ Cell = "E43"; Value = "15.92%" },
    @{ Cell = "D44"; Value = "0.008312" },
    @{ Cell = "E44"; Value = "11.86%" },
    @{ Cell = "D45"; Value = "0.3468" },
    @{ Cell = "E45"; Value = "3.76%" },
    @{ Cell = "D46"; Value = "0.00006705" },
    @{ Cell = "E46"; Value = "-2.42%" },
    @{ Cell = "D47"; Value = "0.00000000753" },
    @{ Cell = "E47"; Value = "0.31%" },
    @{ Cell = "D48"; Value = "0.1958" },
    @{ Cell = "E48"; Value = "43.57%" },
    @{ Cell = "D49"; Value = "0.004217" },
    @{ Cell = "E49"; Value = "40.44%" },
    @{ Cell = "D50"; Value = "0.00002108" },
    @{ Cell = "E50"; Value = "0.31%" },
    @{ Cell = "D51"; Value = "0.0002008" },
    @{ Cell = "E51"; Value = "0.31%" }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = "'" + $u.Value
    $ws.Range($u.Cell).ClearFormats()
}
